# 8 View Ave update
# Adds a "Replacement / mitigation" mini table (rows 18-23, cols E:G) to
# Sheet1, a new currency number format used by the "Mitigation amount"
# row, and updates the current selection / column widths to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New mini "replacement/mitigation" table (E18:G23) -------------------
# Shared-string insertion order matters (it determines the sst indices),
# so write the E-column labels before the F18/G18 headers: that reproduces
# Replacement required(46), Trees planted(47), Mitigation required(48),
# Mitigation amount(49), Today(50), With 100% Replacement(51).
$ws.Range("E20").Value = "Replacement required"
$ws.Range("E21").Value = "Trees planted"
$ws.Range("E22").Value = "Mitigation required"
$ws.Range("E23").Value = "Mitigation amount"
$ws.Range("F18").Value = "Today"
$ws.Range("G18").Value = "With 100% Replacement"

$ws.Range("E19").Value = "Trees removed"
$ws.Range("F19").Formula = "=B46"
$ws.Range("G19").Formula = "=B46"

$ws.Range("F20").Formula = "=ROUND(F19/2, 0)"
$ws.Range("G20").Formula = "=G19"

$ws.Range("F21").Formula = "=J14"
$ws.Range("G21").Formula = "=J14"

$ws.Range("F22").Formula = "=F20-F21"
$ws.Range("G22").Formula = "=G20-G21"

$ws.Range("F23").Formula = "=F22*100"
$ws.Range("G23").Formula = "=G22*100"

# New currency number format ("$"#,##0) for the mitigation-amount row.
$ws.Range("F23:G23").NumberFormat = """$""#,##0"

# --- View state: column widths + selection -------------------------------
$ws.Columns("F").ColumnWidth = 18.666666666666668
$ws.Columns("G").ColumnWidth = 8.333333333333332

$ws.Range("I29").Select()
